$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

# Row 7 (old row 7 -> becomes old row 8's data)
$ws.Range("A7").Value = 111401581
$ws.Range("B7").Value = 90710
$ws.Range("E7").Value = 5449
$ws.Range("F7").Value = "Svart taggsvamp"
$ws.Range("G7").Value = "Phellodon niger"
$ws.Range("H7").Value = "(Fr.:Fr.) P.Karst."
Set-TextValue "I7" "2"
$ws.Range("P7").Value = "Öst Låssbytjärnet, Vrm"
$ws.Range("Q7").Value = 318450.5531044828
$ws.Range("R7").Value = 6596617.106492633

# Row 8 (old row 8 -> becomes old row 9's data)
$ws.Range("A8").Value = 111401607
$ws.Range("B8").Value = 90685
$ws.Range("D8").Value = "VU"
$ws.Range("E8").Value = 1440
$ws.Range("F8").Value = "Brödtaggsvamp"
$ws.Range("G8").Value = "Hydnellum versipelle"
$ws.Range("H8").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
Set-TextValue "I8" "25"
$ws.Range("Q8").Value = 318453.8345372439
$ws.Range("R8").Value = 6596685.151500781

# Row 9 (old row 9 -> becomes old row 7's data)
$ws.Range("A9").Value = 111401477
$ws.Range("B9").Value = 90651
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 1968
$ws.Range("F9").Value = "Grantaggsvamp"
$ws.Range("G9").Value = "Bankera violascens"
$ws.Range("H9").Value = "(Alb. & Schwein. : Fr.) Pouzar"
Set-TextValue "I9" "7"
$ws.Range("P9").Value = "Öster Låssbytjärnet, Vrm"
$ws.Range("Q9").Value = 318589.7492944719
$ws.Range("R9").Value = 6596680.635420629
